# Estado de Cuenta - agrega nuevos trabajadores en mora, actualiza totales
# y desplaza el bloque de firmas hacia abajo.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Inserta 3 filas completas despues de la fila 16 (donde comienzan los
#    datos en blanco), desplazando hacia abajo todo lo que sigue
#    (filas en blanco + bloque de firma con sus celdas combinadas).
$ws.Range("17:19").Insert(-4121)

# 2) Copia el formato (bordes, fuente, relleno, formato numerico) de la
#    ultima fila de datos existente (16) hacia las dos filas nuevas que
#    continuan la tabla en el mismo estilo.
$ws.Range("B16:J16").Copy()
$ws.Range("B17:J17").PasteSpecial(-4122)
$ws.Range("B16:J16").Copy()
$ws.Range("B18:J18").PasteSpecial(-4122)
$ws.Range("B16:J16").Copy()
$ws.Range("B19:J19").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# 3) Contenido de los nuevos trabajadores en mora.
$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "1051447108"
$ws.Range("D17").Value = "CARLOS RUBEN PAJARO MEDINA"
$ws.Range("E17").Value = "2504"
$ws.Range("F17").Value = 56940
$ws.Range("G17").Value = 1423500

$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "8853503"
$ws.Range("D18").Value = "RONALD JOSE MORALES RICARDO"
$ws.Range("E18").Value = "2504"
$ws.Range("F18").Value = 56940
$ws.Range("G18").Value = 1423500

$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = "1192746587"
$ws.Range("D19").Value = "FRANCISCO ANDRES MARTINEZ AGUIRRE"
$ws.Range("E19").Value = "2504"
$ws.Range("F19").Value = 56940
$ws.Range("G19").Value = 1423500

# 4) Totales actualizados (Valor Mora, Cant. Trabajadores, Cant. Periodos).
$ws.Range("E11").Value = 182953
$ws.Range("C13").Value = 4
$ws.Range("F13").Value = 2

# 5) La columna D (Nombre Trabajador) debe ajustarse al nombre mas largo.
$ws.Columns("D").AutoFit()
